# Delete the row containing "b23-24MQ309"/"Create" (row 24) together with
# its blank spacer row (row 25) on the "Sale 22-23" sheet. Excel shifts all
# rows below up by two, and formulas/relative references re-adjust
# automatically, exactly as captured by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sale 22-23")

$ws.Rows("24:25").Delete()

$ws.Activate()
$ws.Range("A27").Select()
